$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 314078
$ws.Cells.Item(2, 4).Value = 400449863
$ws.Cells.Item(4, 3).Value = 312
$ws.Cells.Item(4, 4).Value = 446418
$ws.Cells.Item(8, 3).Value = 844
$ws.Cells.Item(8, 4).Value = 1243607
$ws.Cells.Item(10, 3).Value = 115579
$ws.Cells.Item(10, 4).Value = 169373365
$ws.Cells.Item(12, 3).Value = 58255
$ws.Cells.Item(12, 4).Value = 84083724
$ws.Cells.Item(16, 3).Value = 3924
$ws.Cells.Item(16, 4).Value = 5570261
$ws.Cells.Item(20, 3).Value = 6397
$ws.Cells.Item(20, 4).Value = 8928807
$ws.Cells.Item(22, 3).Value = 76031
$ws.Cells.Item(22, 4).Value = 94920236
$ws.Cells.Item(28, 3).Value = 32096
$ws.Cells.Item(28, 4).Value = 46990193
$ws.Cells.Item(30, 3).Value = 11298
$ws.Cells.Item(30, 4).Value = 16253953
$ws.Cells.Item(33, 3).Value = 1544
$ws.Cells.Item(33, 4).Value = 2167307
$ws.Cells.Item(35, 3).Value = 1750
$ws.Cells.Item(35, 4).Value = 2468260
$ws.Cells.Item(36, 3).Value = 95569
$ws.Cells.Item(36, 4).Value = 120414056
$ws.Cells.Item(38, 3).Value = 82
$ws.Cells.Item(38, 4).Value = 115008
$ws.Cells.Item(44, 3).Value = 43939
$ws.Cells.Item(44, 4).Value = 64406727
$ws.Cells.Item(46, 3).Value = 8987
$ws.Cells.Item(46, 4).Value = 12899920
$ws.Cells.Item(48, 3).Value = 1382
$ws.Cells.Item(48, 4).Value = 1919400
$ws.Cells.Item(51, 3).Value = 2225
$ws.Cells.Item(51, 4).Value = 3102064
$ws.Cells.Item(52, 3).Value = 67752
$ws.Cells.Item(52, 4).Value = 85050831
$ws.Cells.Item(58, 3).Value = 27789
$ws.Cells.Item(58, 4).Value = 40757814
$ws.Cells.Item(61, 3).Value = 10891
$ws.Cells.Item(61, 4).Value = 15746434
$ws.Cells.Item(63, 3).Value = 1341
$ws.Cells.Item(63, 4).Value = 1873789
$ws.Cells.Item(67, 3).Value = 1415
$ws.Cells.Item(67, 4).Value = 1979768
$ws.Cells.Item(69, 3).Value = 20093
$ws.Cells.Item(69, 4).Value = 26320731
$ws.Cells.Item(71, 3).Value = 25
$ws.Cells.Item(71, 4).Value = 36225
$ws.Cells.Item(73, 3).Value = 7458
$ws.Cells.Item(73, 4).Value = 10918861
$ws.Cells.Item(75, 3).Value = 5007
$ws.Cells.Item(75, 4).Value = 7269616
$ws.Cells.Item(78, 3).Value = 138041
$ws.Cells.Item(78, 4).Value = 172238482
$ws.Cells.Item(82, 3).Value = 424
$ws.Cells.Item(82, 4).Value = 619271
$ws.Cells.Item(84, 3).Value = 62780
$ws.Cells.Item(84, 4).Value = 92020516
$ws.Cells.Item(87, 3).Value = 29185
$ws.Cells.Item(87, 4).Value = 42230152
$ws.Cells.Item(89, 3).Value = 2684
$ws.Cells.Item(89, 4).Value = 3865020
$ws.Cells.Item(90, 3).Value = 2714
$ws.Cells.Item(90, 4).Value = 3834205
$ws.Cells.Item(91, 3).Value = 31494
$ws.Cells.Item(91, 4).Value = 42669798
$ws.Cells.Item(95, 3).Value = 7725
$ws.Cells.Item(95, 4).Value = 11360101
$ws.Cells.Item(97, 3).Value = 7000
$ws.Cells.Item(97, 4).Value = 10146480
$ws.Cells.Item(101, 3).Value = 8439
$ws.Cells.Item(101, 4).Value = 11717672
$ws.Cells.Item(103, 3).Value = 2139
$ws.Cells.Item(103, 4).Value = 3150470
$ws.Cells.Item(105, 3).Value = 2882
$ws.Cells.Item(105, 4).Value = 4210919
$ws.Cells.Item(109, 3).Value = 138608
$ws.Cells.Item(109, 4).Value = 171456720
$ws.Cells.Item(115, 3).Value = 52065
$ws.Cells.Item(115, 4).Value = 76336116
$ws.Cells.Item(117, 3).Value = 26443
$ws.Cells.Item(117, 4).Value = 38310315
$ws.Cells.Item(118, 3).Value = 1288
$ws.Cells.Item(118, 4).Value = 1762565
$ws.Cells.Item(121, 3).Value = 2172
$ws.Cells.Item(121, 4).Value = 3049718
$ws.Cells.Item(123, 3).Value = 488640
$ws.Cells.Item(123, 4).Value = 644473878
$ws.Cells.Item(125, 3).Value = 207
$ws.Cells.Item(125, 4).Value = 305496
$ws.Cells.Item(128, 3).Value = 1353
$ws.Cells.Item(128, 4).Value = 2005811
$ws.Cells.Item(130, 3).Value = 203418
$ws.Cells.Item(130, 4).Value = 299044164
$ws.Cells.Item(131, 3).Value = 385
$ws.Cells.Item(131, 4).Value = 574290
$ws.Cells.Item(133, 3).Value = 175532
$ws.Cells.Item(133, 4).Value = 255141397
$ws.Cells.Item(136, 3).Value = 2779
$ws.Cells.Item(136, 4).Value = 3904832
$ws.Cells.Item(138, 3).Value = 6104
$ws.Cells.Item(138, 4).Value = 8622701
$ws.Cells.Item(141, 3).Value = 43450
$ws.Cells.Item(141, 4).Value = 58024573
$ws.Cells.Item(147, 3).Value = 13804
$ws.Cells.Item(147, 4).Value = 20255424
$ws.Cells.Item(148, 3).Value = 3672
$ws.Cells.Item(148, 4).Value = 5296465
$ws.Cells.Item(151, 3).Value = 385
$ws.Cells.Item(151, 4).Value = 553431
$ws.Cells.Item(153, 3).Value = 366
$ws.Cells.Item(153, 4).Value = 515751
$ws.Cells.Item(154, 3).Value = 17061
$ws.Cells.Item(154, 4).Value = 22549202
$ws.Cells.Item(158, 3).Value = 6977
$ws.Cells.Item(158, 4).Value = 10149416
$ws.Cells.Item(160, 3).Value = 4866
$ws.Cells.Item(160, 4).Value = 7003456
$ws.Cells.Item(163, 3).Value = 257
$ws.Cells.Item(163, 4).Value = 368283
$ws.Cells.Item(165, 3).Value = 14457
$ws.Cells.Item(165, 4).Value = 20964860
$ws.Cells.Item(166, 3).Value = 1702
$ws.Cells.Item(166, 4).Value = 2531630
$ws.Cells.Item(167, 3).Value = 231
$ws.Cells.Item(167, 4).Value = 341302
$ws.Cells.Item(171, 3).Value = 86175
$ws.Cells.Item(171, 4).Value = 107828686
$ws.Cells.Item(176, 3).Value = 637
$ws.Cells.Item(176, 4).Value = 938848
$ws.Cells.Item(178, 3).Value = 33451
$ws.Cells.Item(178, 4).Value = 49058394
$ws.Cells.Item(180, 3).Value = 12777
$ws.Cells.Item(180, 4).Value = 18460490
$ws.Cells.Item(182, 3).Value = 1227
$ws.Cells.Item(182, 4).Value = 1716696
$ws.Cells.Item(184, 3).Value = 1589
$ws.Cells.Item(184, 4).Value = 2236693
$ws.Cells.Item(186, 3).Value = 233736
$ws.Cells.Item(186, 4).Value = 290645864
$ws.Cells.Item(188, 3).Value = 164
$ws.Cells.Item(188, 4).Value = 236236
$ws.Cells.Item(192, 3).Value = 864
$ws.Cells.Item(192, 4).Value = 1270997
$ws.Cells.Item(194, 3).Value = 85577
$ws.Cells.Item(194, 4).Value = 125455104
$ws.Cells.Item(197, 3).Value = 32460
$ws.Cells.Item(197, 4).Value = 46718166
$ws.Cells.Item(200, 3).Value = 5010
$ws.Cells.Item(200, 4).Value = 7141248
$ws.Cells.Item(203, 3).Value = 4665
$ws.Cells.Item(203, 4).Value = 6452607
$ws.Cells.Item(206, 3).Value = 258132
$ws.Cells.Item(206, 4).Value = 319543320
$ws.Cells.Item(207, 3).Value = 155
$ws.Cells.Item(207, 4).Value = 169473
$ws.Cells.Item(213, 3).Value = 606
$ws.Cells.Item(213, 4).Value = 882906
$ws.Cells.Item(215, 3).Value = 93817
$ws.Cells.Item(215, 4).Value = 137262984
$ws.Cells.Item(218, 3).Value = 50470
$ws.Cells.Item(218, 4).Value = 72947143
$ws.Cells.Item(221, 3).Value = 4578
$ws.Cells.Item(221, 4).Value = 6424215
$ws.Cells.Item(224, 3).Value = 5503
$ws.Cells.Item(224, 4).Value = 7607738
$ws.Cells.Item(227, 3).Value = 104292
$ws.Cells.Item(227, 4).Value = 130582634
$ws.Cells.Item(234, 3).Value = 48904
$ws.Cells.Item(234, 4).Value = 71652542
$ws.Cells.Item(236, 3).Value = 12156
$ws.Cells.Item(236, 4).Value = 17475692
$ws.Cells.Item(238, 3).Value = 1872
$ws.Cells.Item(238, 4).Value = 2683109
$ws.Cells.Item(240, 3).Value = 2416
$ws.Cells.Item(240, 4).Value = 3375596
$ws.Cells.Item(241, 3).Value = 252184
$ws.Cells.Item(241, 4).Value = 318555526
$ws.Cells.Item(242, 3).Value = 169
$ws.Cells.Item(242, 4).Value = 209459
$ws.Cells.Item(243, 3).Value = 245
$ws.Cells.Item(243, 4).Value = 351957
$ws.Cells.Item(249, 3).Value = 94443
$ws.Cells.Item(249, 4).Value = 138398849
$ws.Cells.Item(250, 3).Value = 77
$ws.Cells.Item(250, 4).Value = 112661
$ws.Cells.Item(252, 3).Value = 63676
$ws.Cells.Item(252, 4).Value = 92279279
$ws.Cells.Item(254, 3).Value = 2361
$ws.Cells.Item(254, 4).Value = 3332247
$ws.Cells.Item(257, 3).Value = 4442
$ws.Cells.Item(257, 4).Value = 6234180
